# Update monitoring RKA 2026 dari web
# Applies the 2026 monitoring workbook data refresh: updates ProgramConsumed (D),
# RealisasiConsumed (E) and RealisasiActual (G) figures across the "Monitoring RKA 2026"
# sheet to reflect the latest pull from the source system, and adjusts the
# active view (selection / scroll) to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monitoring RKA 2026")

# --- Data updates -------------------------------------------------------
    $ws.Range("D2").Value = 50000000
    $ws.Range("D3").Value = 50000000
    $ws.Range("E3").Value = 59570632314
    $ws.Range("D4").Value = 27879534020
    $ws.Range("D5").Value = 150000000
    $ws.Range("D6").Value = 320862500
    $ws.Range("D7").Value = 13426784080
    $ws.Range("D8").Value = 11686187000
    $ws.Range("D9").Value = 25032829880
    $ws.Range("D10").Value = 24914217620
    $ws.Range("D11").Value = 49633981710
    $ws.Range("D12").Value = 33852315000
    $ws.Range("D13").Value = 83354006000
    $ws.Range("D45").Value = 1000000000
    $ws.Range("D56").Value = 500000000
    $ws.Range("D62").Value = 113000000
    $ws.Range("D63").Value = 243000000
    $ws.Range("E63").Value = 171284525
    $ws.Range("G63").Value = 143833500
    $ws.Range("D64").Value = 243000000
    $ws.Range("D65").Value = 1243000000
    $ws.Range("D66").Value = 593000000
    $ws.Range("D67").Value = 1043000000
    $ws.Range("D68").Value = 963000000
    $ws.Range("D69").Value = 443000000
    $ws.Range("D70").Value = 493000000
    $ws.Range("D71").Value = 663000000
    $ws.Range("D72").Value = 293000000
    $ws.Range("D73").Value = 113000000
    $ws.Range("D74").Value = 310000000
    $ws.Range("D75").Value = 310000000
    $ws.Range("G75").Value = 626594334
    $ws.Range("D76").Value = 310000000
    $ws.Range("D77").Value = 420000000
    $ws.Range("D78").Value = 420000000
    $ws.Range("D79").Value = 420000000
    $ws.Range("D80").Value = 420000000
    $ws.Range("D81").Value = 370000000
    $ws.Range("D82").Value = 400000000
    $ws.Range("D83").Value = 400000000
    $ws.Range("D84").Value = 400000000
    $ws.Range("D85").Value = 320000000
    $ws.Range("E87").Value = 10678407138
    $ws.Range("E99").Value = 113020000
    $ws.Range("G99").Value = 34500000
    $ws.Range("E147").Value = 70819316
    $ws.Range("G147").Value = 86687950
    $ws.Range("E159").Value = 222698232
    $ws.Range("G159").Value = 222698232
    $ws.Range("E171").Value = 25211824009
    $ws.Range("G171").Value = 457037427
    $ws.Range("E183").Value = 108905600
    $ws.Range("E231").Value = 37750500
    $ws.Range("G231").Value = 64950500
    $ws.Range("E255").Value = 3188732296
    $ws.Range("G255").Value = 147536656
    $ws.Range("E279").Value = 15775000
    $ws.Range("E315").Value = 44348055
    $ws.Range("G315").Value = 112098409
    $ws.Range("E327").Value = 158708365
    $ws.Range("G327").Value = 158708365
    $ws.Range("E339").Value = 882186785
    $ws.Range("G339").Value = 28044916
    $ws.Range("E351").Value = 28227453
    $ws.Range("E399").Value = 132788748
    $ws.Range("E411").Value = 165926590
    $ws.Range("G411").Value = 165950842
    $ws.Range("E423").Value = 683351141
    $ws.Range("G423").Value = 93332500
    $ws.Range("E483").Value = 82708000
    $ws.Range("G483").Value = 92994741
    $ws.Range("E495").Value = 193379527
    $ws.Range("G495").Value = 193379527
    $ws.Range("E507").Value = 2975215110
    $ws.Range("G507").Value = 258244550
    $ws.Range("E567").Value = 124319935
    $ws.Range("G567").Value = 94911900
    $ws.Range("E579").Value = 197235478
    $ws.Range("G579").Value = 197235478
    $ws.Range("E591").Value = 4701283132
    $ws.Range("G591").Value = 24650000
    $ws.Range("G603").Value = 36500000
    $ws.Range("G615").Value = 62100000
    $ws.Range("E651").Value = 80820599
    $ws.Range("G651").Value = 137491406
    $ws.Range("E663").Value = 158370071
    $ws.Range("G663").Value = 158370071
    $ws.Range("E675").Value = 165605895
    $ws.Range("G675").Value = 86559328
    $ws.Range("E735").Value = 43215908
    $ws.Range("G735").Value = 32679201
    $ws.Range("E747").Value = 141087260
    $ws.Range("G747").Value = 141087260
    $ws.Range("E759").Value = 14723328892
    $ws.Range("G759").Value = 426154218
    $ws.Range("E819").Value = 82038470
    $ws.Range("G819").Value = 70154750
    $ws.Range("E831").Value = 202611220
    $ws.Range("G831").Value = 207808845
    $ws.Range("E843").Value = 316431833
    $ws.Range("G843").Value = 187448104
    $ws.Range("E903").Value = 131493477
    $ws.Range("G903").Value = 44852000
    $ws.Range("E915").Value = 35802840
    $ws.Range("G915").Value = 34478198
    $ws.Range("E927").Value = 660310335
    $ws.Range("E951").Value = 37804500
    $ws.Range("E987").Value = 121455950
    $ws.Range("G987").Value = 91152500
    $ws.Range("E999").Value = 130722605
    $ws.Range("G999").Value = 121222605
    $ws.Range("E1071").Value = 3123150
    $ws.Range("G1071").Value = 3117150
    $ws.Range("E1095").Value = 0
    $ws.Range("G1095").Value = 0
    $ws.Range("E1106").Value = 4033286500
    $ws.Range("G1106").Value = 0
    $ws.Range("E1107").Value = 4210736952
    $ws.Range("G1107").Value = 96556000
    $ws.Range("E1118").Value = 82800000
    $ws.Range("G1118").Value = 0
    $ws.Range("E1119").Value = 45540000
    $ws.Range("G1119").Value = 82800000
    $ws.Range("E1130").Value = 0
    $ws.Range("G1130").Value = 0
    $ws.Range("E1131").Value = 12720000
    $ws.Range("G1131").Value = 0
    $ws.Range("E1143").Value = 0
    $ws.Range("G1143").Value = 0
    $ws.Range("E1154").Value = 142025182
    $ws.Range("G1154").Value = 20219062
    $ws.Range("E1155").Value = 130923900
    $ws.Range("G1155").Value = 163539770
    $ws.Range("E1166").Value = 268881303
    $ws.Range("G1166").Value = 268553262
    $ws.Range("E1167").Value = 186153806
    $ws.Range("G1167").Value = 178077441
    $ws.Range("E1178").Value = 9143456586
    $ws.Range("G1178").Value = 0
    $ws.Range("E1179").Value = 477983925
    $ws.Range("G1179").Value = 166766116
    $ws.Range("E1190").Value = 0
    $ws.Range("G1190").Value = 0
    $ws.Range("E1191").Value = 0
    $ws.Range("G1191").Value = 0
    $ws.Range("E1202").Value = 32200000
    $ws.Range("G1202").Value = 0
    $ws.Range("E1203").Value = 0
    $ws.Range("G1203").Value = 32200000
    $ws.Range("E1214").Value = 18435000
    $ws.Range("G1214").Value = 0
    $ws.Range("E1215").Value = 64735000
    $ws.Range("G1215").Value = 18435000
    $ws.Range("E1227").Value = 0
    $ws.Range("G1227").Value = 0
    $ws.Range("E1238").Value = 78887534
    $ws.Range("G1238").Value = 17317500
    $ws.Range("E1239").Value = 37118745
    $ws.Range("G1239").Value = 97688779
    $ws.Range("E1250").Value = 122102415
    $ws.Range("G1250").Value = 116102415
    $ws.Range("E1251").Value = 108945350
    $ws.Range("G1251").Value = 114601656

# --- View state ----------------------------------------------------------
# Keep the existing frozen-pane split (3 cols / 1 row) but move the active
# selection to where the user left off after refreshing the data.
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("D2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 999
$win.ScrollColumn = 4
$ws.Range("F1005").Select()
